# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column BP ("20-aug") with its 24 hourly values
#  - "Gaz" sheet: append a new row (65) for 2025-08-18
#  - "CO2" sheet: append a new row (65) for 2025-08-18

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column BP
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Header cell BP1 - copy style from the previous header (BO1, bold/centered/
# bordered) so the new column looks like the rest of the header row, then set
# its own text.
$ws1.Range("BO1").Copy($ws1.Range("BP1"))
$ws1.Range("BP1").Value = "20-aug"

# Hourly values for the new day, row 2 through row 25.
$bpValues = @{
    2  = 82.90000000000001
    3  = 74.12
    4  = 67.05
    5  = 66.08
    6  = 64.95
    7  = 70.59999999999999
    8  = 74.06999999999999
    9  = 78.43000000000001
    10 = 89.81
    11 = 71.55
    12 = 49.8
    13 = 56.03
    14 = 49.02
    15 = 35
    16 = 31.53
    17 = 27.2
    18 = 27.99
    19 = 37.72
    20 = 54.11
    21 = 69.27
    22 = 78.73
    23 = 100.07
    24 = 97.52
    25 = 84.98
}

foreach ($row in $bpValues.Keys) {
    $ws1.Cells.Item($row, 68).Value = $bpValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 65 for 2025-08-18
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")
$cellA = $ws2.Cells.Item(65, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-08-18"
$cellA.Style = "Normal"
$ws2.Cells.Item(65, 2).Value = 29.95

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 65 for 2025-08-18
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")
$cellA3 = $ws3.Cells.Item(65, 1)
$cellA3.NumberFormat = "@"
$cellA3.Value = "2025-08-18"
$cellA3.Style = "Normal"
$ws3.Cells.Item(65, 2).Value = 71
